$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values per repull/recalculation of data
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = -9
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = -5
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = -4
$ws.Range("F22").Value = -1
$ws.Range("F23").Value = -10
$ws.Range("F24").Value = 0
$ws.Range("F28").Value = -1
